# Update the [ts_model] column (AJ) values on Sheet1 from 3 to 0 for rows 2-24.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AJ2:AJ24").Value = 0
